# PF7111/hand_faired_curve.xlsx — "added baseline temp plot for comparison"
#
# Updates the hand-faired curve data (cols A/B), extends the table with two
# additional rows, adds a "Slope" column (C) computed from consecutive
# points, and repositions/resizes the scatter chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Updated X (dPp_qcic, col A) / Y (Mic, col B) data, rows 2-15.
#    Rows 14-15 are brand-new data points extending the curve.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 0.291187
$ws.Range("B2").Value = 0.02164

$ws.Range("A3").Value = 0.35
$ws.Range("B3").Value = 0.0172

$ws.Range("A4").Value = 0.4
$ws.Range("B4").Value = 0.0144

$ws.Range("A5").Value = 0.45
$ws.Range("B5").Value = 0.01242

$ws.Range("A6").Value = 0.5
$ws.Range("B6").Value = 0.01097

$ws.Range("A7").Value = 0.55
$ws.Range("B7").Value = 0.00973

$ws.Range("A8").Value = 0.6
$ws.Range("B8").Value = 0.00875

$ws.Range("A9").Value = 0.65
$ws.Range("B9").Value = 0.0078

$ws.Range("A10").Value = 0.7
$ws.Range("B10").Value = 0.00695

$ws.Range("A11").Value = 0.75
$ws.Range("B11").Value = 0.006125

$ws.Range("A12").Value = 0.8
$ws.Range("B12").Value = 0.00535

$ws.Range("A13").Value = 0.85
$ws.Range("B13").Value = 0.00474

$ws.Range("A14").Value = 0.9
$ws.Range("B14").Value = 0.00415

$ws.Range("A15").Value = 0.919487
$ws.Range("B15").Value = 0.00395

# ---------------------------------------------------------------------
# 2. New "Slope" column (C): header + per-row secant slope formula.
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "Slope"

$ws.Range("C3").Formula = "=(B3-B2)/(A3-A2)"
$ws.Range("C4:C15").Formula = "=(B4-B3)/(A4-A3)"

# ---------------------------------------------------------------------
# 3. Move / resize the chart (it now spans a different, slightly
#    smaller range of columns/rows to make room alongside the new data).
# ---------------------------------------------------------------------
$co = $ws.ChartObjects(1)
$co.Left = 236.3125
$co.Top = 7
$co.Width = 962.4375
$co.Height = 559

# ---------------------------------------------------------------------
# 4. Restore the selection the author left the sheet on.
# ---------------------------------------------------------------------
$ws.Range("B3").Select()
